$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Integer serial number" demo rows (A7:A10) ---
# These cells held 0,1,2,46016 and now hold 1,2,46016,2958465 - the old
# largest-serial-number row (11) becomes redundant and is removed below.
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 46016
$ws.Range("A10").Value = 2958465

# Remove the now-duplicate row 11 (it used to hold 2958465 / the "Largest
# accepted serial number" comment, which now belongs on row 10).
$ws.Rows(11).Delete()

# Row 11's comment moves up onto row 10 alongside the value that moved there.
# (D10 previously only carried a style with no content, D11 carried the
# text with the default/no style - clear D10's leftover style to match.)
$ws.Range("D10").Value = "Largest accepted serial number"
$ws.Range("D10").Style = "Normal"

# --- Second table (error examples) ---
# Insert a fresh row after what is now row 17 (DAY(-1)) for the new
# "Zero serial number" example, pushing DAY(2958466) / DAY(DATE(1900,2,29))
# back down to rows 19/20.
$ws.Rows(18).Insert()

# Row 17's example becomes DAY(-5) instead of DAY(-1) (comment stays the same).
$ws.Range("A17").Formula = "=DAY(-5)"
$ws.Range("B17").Formula = "=FORMULATEXT(A17)"

# The newly inserted row 18 documents DAY(0), which no longer errors.
$ws.Range("A18").Formula = "=DAY(0)"
$ws.Range("B18").Formula = "=FORMULATEXT(A18)"
$ws.Range("C18").Value = "Zero serial number"

# --- Selection moved to B12 (now an empty spacer row) ---
$ws.Range("B12").Select()
